$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C20").Value = "40 equal symbols  + </CTR>"
$ws.Range("C20").Select()
